$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cyclic shift of columns J,K,L (Comment moves to the end; SamplePortion/SamplePortionUnit move left)
$ws.Range("J1").Value = 'SamplePortion'
$ws.Range("K1").Value = 'SamplePortionUnit'
$ws.Range("L1").Value = 'Comment'
$ws.Range("J2").Value = '# Prise d''essai'
$ws.Range("K2").Value = '# Unité de mesure de la prise d’essai'
$ws.Range("L2").Value = '# Commentaire'
$ws.Range("J3").Value = '#float'
$ws.Range("K3").Value = '#string'
$ws.Range("L3").Value = '#string'
$ws.Range("J4").Value = '# format: nombre décimal, ne pas spécifier d''unité'
$ws.Range("K4").Value = '# format: texte'
$ws.Range("L4").Value = '# format: texte libre'
$ws.Range("J5").Value = '# ex: 2.5'
$ws.Range("K5").Value = '# ex: mg'
$ws.Range("L5").Value = ""
